# Insert the reference to "Decreto N°140/04 del Ministerio de Salud..." right
# after "...que rigen los actos de los Organos de la Administracion del
# Estado;" and before the "${art8}" merge-field placeholder, in the long
# legal "Vistos" citation paragraph.

$d = $word.ActiveDocument

$anchor = "que rigen los actos de los Órganos de la Administración del Estado;"
$newText = " Decreto N°140/04 del Ministerio de Salud que aprobó el Reglamento orgánico de los Servicios de Salud,"

$rng = $d.Content
$found = $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter($newText)
} else {
    Write-Output "WARNING: anchor text not found; document left unchanged"
}
